# This document was re-saved with a newer XML serializer (Apache POI 3.15)
# during a packaging fix. The upgrade only changed how attributes and
# namespace declarations are ordered when the XML parts are written back
# out - it did not change any visible document content, formatting,
# styles, or structure. There is nothing in the Word object model to
# change here: every paragraph, run, style and section property is
# exactly the same as before the re-save.
$d = $word.ActiveDocument
